$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 1.233604490811928
$ws.Range("C2").Value = 0.0724348790352991
$ws.Range("D2").Value = 0.004012937064468591
$ws.Range("E2").Value = 0.0651270280068843
$ws.Range("F2").Value = 4.639663493202448
$ws.Range("I2").Value = 3.007647047798713
$ws.Range("J2").Value = 0.1738394188916423
$ws.Range("K2").Value = 1.098044654369374
$ws.Range("L2").Value = 0.3540689197456715
# Row 3
$ws.Range("B3").Value = 1.216754902808873
$ws.Range("C3").Value = 0.06715396461659395
$ws.Range("D3").Value = 0.004145253105162894
$ws.Range("E3").Value = 0.06537049419045626
$ws.Range("F3").Value = 4.599691938734992
$ws.Range("I3").Value = 2.983223206827091
$ws.Range("J3").Value = 0.1737394437704296
$ws.Range("K3").Value = 1.075287720699265
$ws.Range("L3").Value = 0.3524579666636924
# Row 4
$ws.Range("B4").Value = 1.207252263149627
$ws.Range("C4").Value = 0.06395599346494407
$ws.Range("D4").Value = 0.004233372198386931
$ws.Range("E4").Value = 0.06554090196970463
$ws.Range("F4").Value = 4.576408819635787
$ws.Range("I4").Value = 2.968916534873571
$ws.Range("J4").Value = 0.1737113378703015
$ws.Range("K4").Value = 1.062085975564287
$ws.Range("L4").Value = 0.3516446564153526
# Row 5
$ws.Range("B5").Value = 1.203592118987217
$ws.Range("C5").Value = 0.06266391640153302
$ws.Range("D5").Value = 0.004271017070007321
$ws.Range("E5").Value = 0.06561561930583792
$ws.Range("F5").Value = 4.567237214676297
$ws.Range("I5").Value = 2.963259630052391
$ws.Range("J5").Value = 0.1737082610830996
$ws.Range("K5").Value = 1.056900152086683
$ws.Range("L5").Value = 0.3513574990491435
# Row 6
$ws.Range("B6").Value = 1.202997182580049
$ws.Range("C6").Value = 0.06245003754631284
$ws.Range("D6").Value = 0.004277372991440487
$ws.Range("E6").Value = 0.06562834507846826
$ws.Range("F6").Value = 4.565733381867048
$ws.Range("I6").Value = 2.962330757970761
$ws.Range("J6").Value = 0.173708256356079
$ws.Range("K6").Value = 1.05605076984196
$ws.Range("L6").Value = 0.3513124926116546
# Row 7
$ws.Range("B7").Value = 1.207202041390559
$ws.Range("C7").Value = 0.06393852312751847
$ws.Range("D7").Value = 0.004233872854511045
$ws.Range("E7").Value = 0.06554188825676999
$ws.Range("F7").Value = 4.576283847233867
$ws.Range("I7").Value = 2.968839543011441
$ws.Range("J7").Value = 0.1737112624480943
$ws.Range("K7").Value = 1.062015252128731
$ws.Range("L7").Value = 0.3516406043677449
# Row 8
$ws.Range("B8").Value = 1.227619916911777
$ws.Range("C8").Value = 0.07060472970748322
$ws.Range("D8").Value = 0.004057136276964801
$ws.Range("E8").Value = 0.06520664220992423
$ws.Range("F8").Value = 4.62561969063205
$ws.Range("I8").Value = 2.9990824109492
$ws.Range("J8").Value = 0.1737980450871355
$ws.Range("K8").Value = 1.090038064358481
$ws.Range("L8").Value = 0.3534770064475481
# Row 9
$ws.Range("B9").Value = 1.274341815827427
$ws.Range("C9").Value = 0.08403464828960239
$ws.Range("D9").Value = 0.003764836722543574
$ws.Range("E9").Value = 0.0647145845617132
$ws.Range("F9").Value = 4.732381880448486
$ws.Range("I9").Value = 3.06387632622139
$ws.Range("J9").Value = 0.1742320158726756
$ws.Range("K9").Value = 1.151110071799508
$ws.Range("L9").Value = 0.3584712324410049
# Row 10
$ws.Range("B10").Value = 1.312739641244974
$ws.Range("C10").Value = 0.09412641534478894
$ws.Range("D10").Value = 0.003582800054322988
$ws.Range("E10").Value = 0.06445307617084417
$ws.Range("F10").Value = 4.816965553773656
$ws.Range("I10").Value = 3.114854153546688
$ws.Range("J10").Value = 0.1747114558828642
$ws.Range("K10").Value = 1.199719168271145
$ws.Range("L10").Value = 0.3629880283306619
# Row 11
$ws.Range("B11").Value = 1.331091868551681
$ws.Range("C11").Value = 0.09876778987953116
$ws.Range("D11").Value = 0.00350701291975497
$ws.Range("E11").Value = 0.06435565919582764
$ws.Range("F11").Value = 4.856788753741
$ws.Range("I11").Value = 3.138784384221537
$ws.Range("J11").Value = 0.1749644055659552
$ws.Range("K11").Value = 1.222647243526296
$ws.Range("L11").Value = 0.3652265989876895
# Row 12
$ws.Range("B12").Value = 1.338168507203818
$ws.Range("C12").Value = 0.10053273463339
$ws.Range("D12").Value = 0.003479317905806134
$ws.Range("E12").Value = 0.06432185417751057
$ws.Range("F12").Value = 4.872062799142924
$ws.Range("I12").Value = 3.147952963583137
$ws.Range("J12").Value = 0.1750651971898947
$ws.Range("K12").Value = 1.231446856668612
$ws.Range("L12").Value = 0.3661006867061189
# Row 13
$ws.Range("B13").Value = 1.336638778933803
$ws.Range("C13").Value = 0.1001522936563788
$ws.Range("D13").Value = 0.003485237964576005
$ws.Range("E13").Value = 0.06432899771318823
$ws.Range("F13").Value = 4.868764632249963
$ws.Range("I13").Value = 3.145973594973142
$ws.Range("J13").Value = 0.175043267407915
$ws.Range("K13").Value = 1.229546488205074
$ws.Range("L13").Value = 0.3659112633562813
# Row 14
$ws.Range("B14").Value = 1.331671522712782
$ws.Range("C14").Value = 0.09891284501978248
$ws.Range("D14").Value = 0.003504714348553151
$ws.Range("E14").Value = 0.0643528162914464
$ws.Range("F14").Value = 4.858041471253927
$ws.Range("I14").Value = 3.139536548283857
$ws.Range("J14").Value = 0.1749725974804903
$ws.Range("K14").Value = 1.223368843281918
$ws.Range("L14").Value = 0.3652979820988236
# Row 15
$ws.Range("B15").Value = 1.328645474702711
$ws.Range("C15").Value = 0.09815460753557659
$ws.Range("D15").Value = 0.003516774763804076
$ws.Range("E15").Value = 0.06436780717881874
$ws.Range("F15").Value = 4.851498483060141
$ws.Range("I15").Value = 3.13560757798362
$ws.Range("J15").Value = 0.1749299617313937
$ws.Range("K15").Value = 1.219600125971539
$ws.Range("L15").Value = 0.3649257648967676
# Row 16
$ws.Range("B16").Value = 1.311558068985136
$ws.Range("C16").Value = 0.09382411552010694
$ws.Range("D16").Value = 0.003587893724356306
$ws.Range("E16").Value = 0.06445987498163674
$ws.Range("F16").Value = 4.814390122761097
$ws.Range("I16").Value = 3.113305175039059
$ws.Range("J16").Value = 0.1746956257548788
$ws.Range("K16").Value = 1.198237173153331
$ws.Range("L16").Value = 0.3628454279467519
# Row 17
$ws.Range("B17").Value = 1.301301993251656
$ws.Range("C17").Value = 0.09118051272585603
$ws.Range("D17").Value = 0.003633317270850256
$ws.Range("E17").Value = 0.06452186408599303
$ws.Range("F17").Value = 4.791970256095965
$ws.Range("I17").Value = 3.09981314918393
$ws.Range("J17").Value = 0.1745607896124604
$ws.Range("K17").Value = 1.185340554464716
$ws.Range("L17").Value = 0.3616162637167406
# Row 18
$ws.Range("B18").Value = 1.295486268410627
$ws.Range("C18").Value = 0.08966473799705454
$ws.Range("D18").Value = 0.003660105008436165
$ws.Range("E18").Value = 0.06455954724294521
$ws.Range("F18").Value = 4.779201580168092
$ws.Range("I18").Value = 3.092122543051616
$ws.Range("J18").Value = 0.1744865160506812
$ws.Range("K18").Value = 1.177999525345115
$ws.Range("L18").Value = 0.3609265893473008
# Row 19
$ws.Range("B19").Value = 1.293531476572667
$ws.Range("C19").Value = 0.08915233697518943
$ws.Range("D19").Value = 0.003669288636345858
$ws.Range("E19").Value = 0.06457265494759667
$ws.Range("F19").Value = 4.774900065283759
$ws.Range("I19").Value = 3.089530589768003
$ws.Range("J19").Value = 0.1744619319429326
$ws.Range("K19").Value = 1.175527164412358
$ws.Range("L19").Value = 0.3606960516987812
# Row 20
$ws.Range("B20").Value = 1.302385150372288
$ws.Range("C20").Value = 0.09146143595616252
$ws.Range("D20").Value = 0.003628413456924218
$ws.Range("E20").Value = 0.06451505537860669
$ws.Range("F20").Value = 4.79434378004288
$ws.Range("I20").Value = 3.101242187791541
$ws.Range("J20").Value = 0.1745748036650916
$ws.Range("K20").Value = 1.186705476903654
$ws.Range("L20").Value = 0.3617453193554354
# Row 21
$ws.Range("B21").Value = 1.333127079716718
$ws.Range("C21").Value = 0.09927670079161999
$ws.Range("D21").Value = 0.003498966462851971
$ws.Range("E21").Value = 0.06434573658644105
$ws.Range("F21").Value = 4.861185859273036
$ws.Range("I21").Value = 3.141424366282379
$ws.Range("J21").Value = 0.1749932191649961
$ws.Range("K21").Value = 1.225180185421664
$ws.Range("L21").Value = 0.3654774018473432
# Row 22
$ws.Range("B22").Value = 1.353959159478336
$ws.Range("C22").Value = 0.1044273496275139
$ws.Range("D22").Value = 0.003420214365914109
$ws.Range("E22").Value = 0.06425305069650022
$ws.Range("F22").Value = 4.906001179382827
$ws.Range("I22").Value = 3.168307924928314
$ws.Range("J22").Value = 0.1752958484027971
$ws.Range("K22").Value = 1.251009040047506
$ws.Range("L22").Value = 0.3680703305262938
# Row 23
$ws.Range("B23").Value = 1.342772996903363
$ws.Range("C23").Value = 0.1016743961668567
$ws.Range("D23").Value = 0.003461712573174847
$ws.Range("E23").Value = 0.06430087869503431
$ws.Range("F23").Value = 4.881978870139505
$ws.Range("I23").Value = 3.153902638996385
$ws.Range("J23").Value = 0.1751316624360655
$ws.Range("K23").Value = 1.237161168987456
$ws.Range("L23").Value = 0.3666723790550463
# Row 24
$ws.Range("B24").Value = 1.301895203872931
$ws.Range("C24").Value = 0.09133441790045538
$ws.Range("D24").Value = 0.003630628374850176
$ws.Range("E24").Value = 0.06451812722509942
$ws.Range("F24").Value = 4.793270333592574
$ws.Range("I24").Value = 3.100595913434674
$ws.Range("J24").Value = 0.174568457802458
$ws.Range("K24").Value = 1.186088166809583
$ws.Range("L24").Value = 0.3616869203811177
# Row 25
$ws.Range("B25").Value = 1.260987195016241
$ws.Range("C25").Value = 0.08036252016991341
$ws.Range("D25").Value = 0.003838139503869353
$ws.Range("E25").Value = 0.06483007966767929
$ws.Range("F25").Value = 4.702423590696071
$ws.Range("I25").Value = 3.045757439977621
$ws.Range("J25").Value = 0.1740863970478621
$ws.Range("K25").Value = 1.133932460599681
$ws.Range("L25").Value = 0.3569712075158833
